$wb = $excel.ActiveWorkbook

# Add the new worksheet "Estimation_LL" after the existing sheet
$ws2 = $wb.Worksheets.Add()
$ws2.Name = "Estimation_LL"

$ws2.Range("A1").Value = "Complexity"
$ws2.Range("B1").Value = "Weights"
$ws2.Range("C1").Value = "UT"
$ws2.Range("D1").Value = "Table"
$ws2.Range("E1").Value = "Valid"
$ws2.Range("F1").Value = "Trans"

$ws2.Range("A2").Value = "Simple "
$ws2.Range("B2").Value = 1
$ws2.Range("C2").Value = 1
$ws2.Range("D2").Value = 0.5
$ws2.Range("E2").Value = 0.25
$ws2.Range("F2").Value = 5

$ws2.Range("A3").Value = "Medium"
$ws2.Range("B3").Value = 1.2
$ws2.Range("C3").Value = 2
$ws2.Range("D3").Value = 1
$ws2.Range("E3").Value = 0.5
$ws2.Range("F3").Value = 10

$ws2.Range("A4").Value = "Complex"
$ws2.Range("B4").Value = 1.5
$ws2.Range("C4").Value = 3
$ws2.Range("D4").Value = 1.5
$ws2.Range("E4").Value = 1
$ws2.Range("F4").Value = 15
